$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final dataset (account, name, balance) in the order it should appear after the edit.
$data = @(
    @("005305448","ALPHASITIO",365703.31),
    @("004479287","ANA",65727.36),
    @("004224011","THOMAS",54291.07),
    @("005103059","WALQUIRIA",38809.89),
    @("004454365","RAFAEL",13735.23),
    @("004452597","LARA",12966.82),
    @("005046790","BEATRIZ",10469.72),
    @("004643153","CARLA",8581.46),
    @("004361159","HFR",5714.31),
    @("005381719","MARIA",5035.12),
    @("004212438","KENIA",3750),
    @("004332747","LOHRAN",2999.22),
    @("004870019","MARIA",1158.49),
    @("004488571","CARLOS",1000),
    @("004240014","ISABELE",998.7),
    @("004334158","LEONE",994.66),
    @("001882235","LAGO",966.83),
    @("004457389","RAFAEL",926.76),
    @("004487016","ROGERIO",921.71),
    @("004392159","RODRIGO",900.21),
    @("004855570","LUISA",895.14),
    @("004550605","REJANE",869.47),
    @("004936634","LEONARDO",865.79),
    @("004211368","ILTON",856.39),
    @("004216504","WANDER",850.24),
    @("004563252","FERNANDO",837.81),
    @("005245032","ROSA",824.46),
    @("002064834","RAFAELA",813.54),
    @("004452476","IVONE",768.76),
    @("004381180","HFR",743.31),
    @("004388077","WLADMIR",703.17),
    @("004359408","HEPTA",693.95),
    @("004646727","RENATA",679.08),
    @("004975924","SERGIO",672.05),
    @("004972070","MARIA",670.51),
    @("004252768","ALESSANDRO",656.22),
    @("004481463","MARA",637.03),
    @("004228456","FLASH",611.4),
    @("004517080","TATIANA",607.94),
    @("005079311","JOVINO",603.35),
    @("005046919","MARIANA",600),
    @("004574428","GUILHERME",596.31),
    @("004806244","CARLA",585.78),
    @("005142592","ALBERTO",551.34),
    @("004238436","DIEGO",547.62),
    @("004395314","MARIA",522.36),
    @("004398174","DANIELE",496.52),
    @("004322719","GISELA",493.87),
    @("004448303","NASSIM",488.3),
    @("004387250","MONICA",478),
    @("005338054","ELAINE",475.92),
    @("000772433","MARCELO",465.56),
    @("005266369","EG",459.39),
    @("005203562","ROBERIO",453.65),
    @("004893911","FAUSTO",449.88),
    @("004415557","FILIPE",427.15),
    @("004276856","DAURA",426.61),
    @("004556150","MARINA",409.35),
    @("004480970","ALBERTO",399.15),
    @("004756981","MATEUS",357.49),
    @("004587511","CARLOS",352.61),
    @("004289402","LARISSA",349.86),
    @("004212476","MARIA",349.25),
    @("004204500","EDWARD",343.71),
    @("005009992","ALINE",330.17),
    @("005101676","ELENI",314.83),
    @("004214592","MERG",312.75),
    @("004480134","JOSE",308.81),
    @("004332103","JOSE",300.31),
    @("005121919","JORGE",297.95),
    @("004259659","BENTO",293.95),
    @("000330949","RENATO",285.35),
    @("000626491","FELIPE",280.06),
    @("004355790","MINEIA",279.74),
    @("004486497","ELENA",257.45),
    @("004927044","CINTIA",257.37),
    @("004870976","HFR",222.74),
    @("004580355","LARISSA",210.77),
    @("004526450","MSD",205.85),
    @("004334062","MERG",205.56),
    @("004475395","DAVID",185.02),
    @("004328934","VALERIA",182.14),
    @("004360431","CARLOS",164.01),
    @("004511696","KRYSCIA",150.47),
    @("005022526","ALEXANDRE",147.18),
    @("004493324","DANIEL",143.45),
    @("005141215","KARINA",137.66),
    @("005274028","RAFAEL",132.43),
    @("004435987","MARCO",125.33),
    @("004211911","ZENILDA",120),
    @("004404342","ADSON",115.85),
    @("004754920","LUIS",114.69),
    @("004421636","PATRICIA",110),
    @("004221638","CAROLINE",109.24),
    @("004536602","TATIANY",108.62),
    @("005040864","ANDRE",100),
    @("004908680","ELENE",99.31),
    @("004472076","RUBENS",99.18),
    @("004339183","JALISON",95.69),
    @("004431591","MARIO",93.87),
    @("005256849","SANDRO",92.78),
    @("004335031","EDMUNDO",92.73),
    @("004350197","GISELA",91.94),
    @("004239387","LUIZ",89.82),
    @("004212132","JOAO",86.38),
    @("004207374","ANGELICA",85.13),
    @("005035754","JOSE",83.31),
    @("004216657","JOAO",80.63),
    @("004318604","RENAN",80.51),
    @("004451996","ADRIANO",80.36),
    @("004267976","E3",79.84),
    @("001294033","VIVIANE",79.82),
    @("005009922","ANA",79.02),
    @("005070742","JUSCELINO",78.29),
    @("004470679","RODOLFO",77.51),
    @("003115072","VICTOR",69.1),
    @("005133039","PAULO",66.51),
    @("004855596","MARIANA",64.36),
    @("004335251","EDMUNDO",62.39),
    @("000834301","MARCUS",57.13),
    @("004588677","RACHEL",55.91),
    @("004517506","LUIZ",55.87),
    @("004215217","CAROLINA",55.66),
    @("004321092","DANIEL",55.23),
    @("004329229","GABRIEL",54.91),
    @("005032151","ANA",52.9),
    @("004268684","PATRICIA",52.7),
    @("004400640","FELIPE",51.44),
    @("004115403","HEBERT",50.87),
    @("004208447","LEILA",50),
    @("004278033","DAISY",47.37),
    @("005216881","RENAN",46.76),
    @("001719494","LUIS",46.74),
    @("004277637","LARA",46.2),
    @("001731007","GUILHERME",44.59),
    @("004581652","CINCO",44.13),
    @("004805133","PATRICIA",41.48),
    @("004958578","ASSAKO",40.7),
    @("005165116","ANA",40.11),
    @("004998717","GIOVANE",40.08),
    @("004238164","DANIELA",38.3),
    @("004520100","ALEXANDRE",37.96),
    @("002401479","JULIO",37.84),
    @("004001621","DANIELA",37.58),
    @("005000656","LUCIA",35.88),
    @("004211922","CARLOS",34.71),
    @("004340984","RENATA",34),
    @("004994036","BALTASAR",33.73),
    @("004691225","ANNA",33.64),
    @("004472431","LUIS",33.08),
    @("004213943","ELISA",30.26),
    @("005018038","ELAINE",27.91),
    @("005105172","VALDIVINO",26.42),
    @("004377415","ANGELA",26.37),
    @("004240292","MARCO",24.3),
    @("004404724","LEANDRO",24.14),
    @("005173958","VENIA",23.27),
    @("004452507","DANIELA",22.74),
    @("005186167","ANDREA",21),
    @("004214604","MARIA",20.75),
    @("004920447","MARILIA",20.67),
    @("004458604","FABIOLA",20),
    @("004204255","AMADO",18.77),
    @("004368994","CRISTINA",18.56),
    @("004756968","DANIELY",18.08),
    @("004381194","ALINNE",17.71),
    @("001879977","THAISSA",17.14),
    @("005143579","GABRIEL",16.18),
    @("005169333","EDUARDO",16.12),
    @("004422594","WANDIR",14.67),
    @("004565108","GUSTAVO",14.51),
    @("000827730","LUCIANA",13.29),
    @("004752461","SERGIO",10.77),
    @("004216298","FLORDELIZ",9.75),
    @("004527606","MARCIA",9.52),
    @("004264780","MARCELO",8.99),
    @("004921978","ELAINE",8.08),
    @("004381415","JOAO",7.95),
    @("004214460","MARIA",7.54),
    @("004693631","NELY",7.36),
    @("004530494","ROSANGELA",6.94),
    @("004854496","JOSE",6.64),
    @("004448501","JOAO",5.55),
    @("004216434","JAIME",4.83),
    @("005142624","RODRIGO",4.75),
    @("004239624","NINA",4.29),
    @("004848927","ULDARICO",3.62),
    @("005142661","SABRINA",3.6),
    @("004322549","SIMONE",3.42),
    @("004382374","THEOMAR",1.98),
    @("005341184","BRENO",1.85),
    @("004332783","IRON",1.73),
    @("004886366","RENATO",1.57),
    @("004308815","ZELI",1.25),
    @("005228239","DEBORA",0.85),
    @("004223502","BRUNA",0.78),
    @("004229526","EDUARDA",0.74),
    @("004212581","MARIA",0.59),
    @("002687737","JOSE",0.5),
    @("004453302","ISABELLA",0.39),
    @("004806286","VERA",0.19),
    @("004371857","NAZARETH",0.18),
    @("004332207","IRACY",0.16),
    @("004357159","JOAO",0.15),
    @("004320840","NATALIA",0.14),
    @("004466350","RAQUEL",0.11),
    @("005047946","GABRIEL",0.09),
    @("004589311","CLARICE",0.06),
    @("004321016","JOAQUIM",0.02),
    @("004850070","RENATO",0.02),
    @("002878817","GUILHERME",0.01),
    @("004400000","VILMA",0.01),
    @("004459875","HELVECIO",0.01),
    @("004612043","YURI",0.01)
)

# Make room: the sheet currently has 208 data rows (rows 2-209), a blank separator
# row (210) and a footer row (211); the new dataset has 214 rows, so insert
# the extra rows needed right before the blank separator row.
$oldDataRows = 208
$newDataRows = $data.Length
$extra = $newDataRows - $oldDataRows
if ($extra -gt 0) {
    $insertAt = 2 + $oldDataRows
    $ws.Range("A" + $insertAt + ":A" + ($insertAt + $extra - 1)).EntireRow.Insert()
}

# Write every data row, forcing column A to text so the zero-padded account
# numbers survive.
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = 2 + $i
    $item = $data[$i]
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $item[0]
    $ws.Cells.Item($r, 2).Value = $item[1]
    $ws.Cells.Item($r, 3).Value = $item[2]
}

# Blank separator + filter-summary footer row, right after the data.
$blankRow = 2 + $data.Length
$footerRow = $blankRow + 1
$ws.Cells.Item($footerRow, 1).Value = "Filtros aplicados:`nDataFim é (Em branco)`nnr_saldo_disponivel não é 0`nPosição é Posição D-1`nDataFim é (Em branco)`nCARTEIRA não está em branco`nNR_CONTA não está em branco`nTIPO_LANCAMENTO não é ED, ET ou Liquidação Doador"
